$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 857
$ws.Range("J12").Value = 749.5
$ws.Range("L12").Value = 749.5
$ws.Range("N12").Value = -1089.5
# Row 17
$ws.Range("H17").Value = 1030.7858
$ws.Range("J17").Value = 1030.7858
$ws.Range("L17").Value = 3092.3574
$ws.Range("N17").Value = -3428.3574
# Row 21
$ws.Range("H21").Value = 16444.334
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -532
# Row 23
$ws.Range("H23").Value = 16444.334
$ws.Range("I23").Value = 1000
$ws.Range("K23").Value = 1000
$ws.Range("M23").Value = -766
# Row 33
$ws.Range("H33").Value = 1562.5834
$ws.Range("I33").Value = 205.2
$ws.Range("J33").Value = 2532.1428
$ws.Range("K33").Value = 205.2
$ws.Range("L33").Value = 2532.1428
$ws.Range("M33").Value = 23.80000000000001
$ws.Range("N33").Value = -2990.1428
# Row 40
$ws.Range("H40").Value = 3625.182
$ws.Range("J40").Value = 3872.125
$ws.Range("L40").Value = 3872.125
$ws.Range("N40").Value = -4222.125
# Row 70
$ws.Range("H70").Value = 12505034
$ws.Range("J70").Value = 8444.777
$ws.Range("L70").Value = 25334.331
$ws.Range("N70").Value = -25874.331
# Row 73
$ws.Range("H73").Value = 12505034
$ws.Range("J73").Value = 8444.777
$ws.Range("L73").Value = 25334.331
$ws.Range("N73").Value = -27206.331
# Row 138
$ws.Range("H138").Value = 1110.8485
$ws.Range("I138").Value = 1040.0938
$ws.Range("K138").Value = 3120.2814
$ws.Range("M138").Value = 2019.7186
# Row 141
$ws.Range("H141").Value = 4252.1875
$ws.Range("I141").Value = 4004.625
$ws.Range("K141").Value = 12013.875
$ws.Range("M141").Value = -6833.875

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1398.6
$ws.Range("I2").Value = 653.3333
$ws.Range("J2").Value = 2516.5
$ws.Range("K2").Value = 653.3333
$ws.Range("L2").Value = 2516.5
$ws.Range("M2").Value = -540.3333
$ws.Range("N2").Value = -2742.5
# Row 32
$ws.Range("H32").Value = 1873.8948
$ws.Range("I32").Value = 1873.8948
$ws.Range("K32").Value = 1873.8948
$ws.Range("M32").Value = -1586.8948
# Row 45
$ws.Range("H45").Value = 1409.7273
$ws.Range("I45").Value = 1065.8889
$ws.Range("J45").Value = 2957
$ws.Range("K45").Value = 1065.8889
$ws.Range("L45").Value = 2957
$ws.Range("M45").Value = -688.8888999999999
$ws.Range("N45").Value = -3711
# Row 61
$ws.Range("H61").Value = 4061.2
$ws.Range("I61").Value = 4061.2
$ws.Range("K61").Value = 4061.2
$ws.Range("M61").Value = -3849.2
# Row 63
$ws.Range("H63").Value = 4125.2173
$ws.Range("I63").Value = 1431.8462
$ws.Range("J63").Value = 7626.6
$ws.Range("K63").Value = 1431.8462
$ws.Range("L63").Value = 7626.6
$ws.Range("M63").Value = -745.8462
$ws.Range("N63").Value = -8998.6
# Row 66
$ws.Range("H66").Value = 4125.2173
$ws.Range("I66").Value = 1431.8462
$ws.Range("J66").Value = 7626.6
$ws.Range("K66").Value = 7159.231
$ws.Range("L66").Value = 38133
$ws.Range("M66").Value = -3727.231
$ws.Range("N66").Value = -44997
# Row 102
$ws.Range("H102").Value = 3787.1155
$ws.Range("I102").Value = 2279
$ws.Range("K102").Value = 2279
$ws.Range("M102").Value = -657
# Row 116
$ws.Range("H116").Value = 1398.6
$ws.Range("I116").Value = 653.3333
$ws.Range("J116").Value = 2516.5
$ws.Range("K116").Value = 653.3333
$ws.Range("L116").Value = 2516.5
$ws.Range("M116").Value = 1640.6667
$ws.Range("N116").Value = -7104.5
# Row 132
$ws.Range("H132").Value = 1617.6229
$ws.Range("I132").Value = 1629.8182
$ws.Range("J132").Value = 1505.8334
$ws.Range("K132").Value = 4889.4546
$ws.Range("L132").Value = 4517.5002
$ws.Range("M132").Value = -2359.4546
$ws.Range("N132").Value = -9577.5002
# Row 136
$ws.Range("H136").Value = 4061.2
$ws.Range("I136").Value = 4061.2
$ws.Range("K136").Value = 12183.6
$ws.Range("M136").Value = -9633.599999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1398.6
$ws.Range("I3").Value = 653.3333
$ws.Range("J3").Value = 2516.5
$ws.Range("K3").Value = 653.3333
$ws.Range("L3").Value = 2516.5
$ws.Range("M3").Value = -539.3333
$ws.Range("N3").Value = -2744.5
# Row 20
$ws.Range("H20").Value = 961.6667
$ws.Range("I20").Value = 985.35297
$ws.Range("K20").Value = 985.35297
$ws.Range("M20").Value = -738.35297
# Row 82
$ws.Range("H82").Value = 16376.444
$ws.Range("I82").Value = 10923.5
$ws.Range("K82").Value = 10923.5
$ws.Range("M82").Value = -10540.5
# Row 85
$ws.Range("H85").Value = 16376.444
$ws.Range("I85").Value = 10923.5
$ws.Range("K85").Value = 10923.5
$ws.Range("M85").Value = -9597.5
# Row 94
$ws.Range("H94").Value = 3844
$ws.Range("I94").Value = 2964.7144
$ws.Range("J94").Value = 9999
$ws.Range("K94").Value = 2964.7144
$ws.Range("L94").Value = 9999
$ws.Range("M94").Value = -2513.7144
$ws.Range("N94").Value = -10901
# Row 96
$ws.Range("H96").Value = 37138
$ws.Range("I96").Value = 13209.5
$ws.Range("K96").Value = 13209.5
$ws.Range("M96").Value = -10463.5
# Row 107
$ws.Range("H107").Value = 4813.857
$ws.Range("I107").Value = 2551
$ws.Range("K107").Value = 2551
$ws.Range("M107").Value = -631

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3201.5
$ws.Range("I58").Value = 2734.5
$ws.Range("J58").Value = 4135.5
$ws.Range("K58").Value = 2734.5
$ws.Range("L58").Value = 4135.5
$ws.Range("M58").Value = -2531.5
$ws.Range("N58").Value = -4541.5
# Row 60
$ws.Range("H60").Value = 34985.125
$ws.Range("J60").Value = 45498.75
$ws.Range("L60").Value = 45498.75
$ws.Range("N60").Value = -46520.75
# Row 132
$ws.Range("H132").Value = 1857.5294
$ws.Range("I132").Value = 1613.0714
$ws.Range("K132").Value = 4839.2142
$ws.Range("M132").Value = -2309.2142
# Row 136
$ws.Range("H136").Value = 3201.5
$ws.Range("I136").Value = 2734.5
$ws.Range("J136").Value = 4135.5
$ws.Range("K136").Value = 8203.5
$ws.Range("L136").Value = 12406.5
$ws.Range("M136").Value = -5653.5
$ws.Range("N136").Value = -17506.5

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 335.16666
$ws.Range("I14").Value = 335.16666
$ws.Range("K14").Value = 1005.49998
$ws.Range("M14").Value = -832.4999799999999
# Row 23
$ws.Range("H23").Value = 3727.4546
$ws.Range("I23").Value = 6711.3335
$ws.Range("J23").Value = 2608.5
$ws.Range("K23").Value = 20134.0005
$ws.Range("L23").Value = 7825.5
$ws.Range("M23").Value = -19899.0005
$ws.Range("N23").Value = -8295.5
# Row 40
$ws.Range("H40").Value = 63.875
$ws.Range("I40").Value = 46.5
$ws.Range("J40").Value = 81.25
$ws.Range("K40").Value = 186
$ws.Range("L40").Value = 325
$ws.Range("M40").Value = -117
$ws.Range("N40").Value = -463
# Row 46
$ws.Range("H46").Value = 5981398
$ws.Range("J46").Value = 1540231.1
$ws.Range("L46").Value = 4620693.300000001
$ws.Range("N46").Value = -4620875.300000001
# Row 87
$ws.Range("H87").Value = 3543
$ws.Range("I87").Value = 3543
$ws.Range("K87").Value = 10629
$ws.Range("M87").Value = -9381
# Row 90
$ws.Range("H90").Value = 3543
$ws.Range("I90").Value = 3543
$ws.Range("K90").Value = 31887
$ws.Range("M90").Value = -25647
# Row 116
$ws.Range("H116").Value = 724.5
$ws.Range("I116").Value = 724.5
$ws.Range("K116").Value = 2173.5
$ws.Range("M116").Value = 1268.5
# Row 118
$ws.Range("H118").Value = 6833.3335
$ws.Range("I118").Value = 4500
$ws.Range("J118").Value = 11500
$ws.Range("K118").Value = 13500
$ws.Range("L118").Value = 34500
$ws.Range("M118").Value = -12257
$ws.Range("N118").Value = -36986
# Row 122
$ws.Range("H122").Value = 2906.923
$ws.Range("J122").Value = 3057.5
$ws.Range("L122").Value = 27517.5
$ws.Range("N122").Value = -32417.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 66672390
$ws.Range("I70").Value = 83338590
$ws.Range("K70").Value = 83338590
$ws.Range("M70").Value = -83338320
# Row 73
$ws.Range("H73").Value = 66672390
$ws.Range("I73").Value = 83338590
$ws.Range("K73").Value = 83338590
$ws.Range("M73").Value = -83337654
# Row 97
$ws.Range("H97").Value = 13195.25
$ws.Range("J97").Value = 25577.334
$ws.Range("L97").Value = 25577.334
$ws.Range("N97").Value = -26569.334

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3644.647
$ws.Range("I7").Value = 3011.923
$ws.Range("J7").Value = 5701
$ws.Range("K7").Value = 3011.923
$ws.Range("L7").Value = 5701
$ws.Range("M7").Value = -2899.923
$ws.Range("N7").Value = -5925
# Row 32
$ws.Range("H32").Value = 15782.5
$ws.Range("I32").Value = 20000
$ws.Range("K32").Value = 20000
$ws.Range("M32").Value = -19683
# Row 55
$ws.Range("H55").Value = 2356.182
$ws.Range("I55").Value = 356.25
$ws.Range("K55").Value = 356.25
$ws.Range("M55").Value = -183.25
# Row 122
$ws.Range("H122").Value = 7196.857
$ws.Range("I122").Value = 6475.6
$ws.Range("K122").Value = 19426.8
$ws.Range("M122").Value = -16976.8
# Row 126
$ws.Range("H126").Value = 3644.647
$ws.Range("I126").Value = 3011.923
$ws.Range("J126").Value = 5701
$ws.Range("K126").Value = 9035.769
$ws.Range("L126").Value = 17103
$ws.Range("M126").Value = -6565.769
$ws.Range("N126").Value = -22043
# Row 132
$ws.Range("H132").Value = 3512.389
$ws.Range("I132").Value = 3322.8333
$ws.Range("J132").Value = 3891.5
$ws.Range("K132").Value = 9968.499899999999
$ws.Range("L132").Value = 11674.5
$ws.Range("M132").Value = -7438.499899999999
$ws.Range("N132").Value = -16734.5

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 2519.4707
$ws.Range("I107").Value = 2497.6428
$ws.Range("J107").Value = 2621.3333
$ws.Range("K107").Value = 7492.928400000001
$ws.Range("L107").Value = 7863.999899999999
$ws.Range("M107").Value = -5572.928400000001
$ws.Range("N107").Value = -11703.9999
